# "Uploaded the revised Event table"
# - Removes the obsolete "sends an update of escalated ticket details" event row
# - Rewords a few cells describing the ticket-escalation transform / reason entry
# - Adds banded row shading (header / alternating rows / highlighted MODULE row)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# xlCenter = -4108, xlLeft = -4131, xlTop = -4160, xlBottom = -4107
$xlCenter = -4108
$xlLeft = -4131

# Color constants (BGR-encoded OLE values) matching the target palette:
#   Header     -> "Blue, Accent 5, Lighter 80%"  (#DEEBF7)
#   Banded row -> "White, Background 1, Darker 5%" (#F2F2F2)
#   Highlight  -> custom lavender (#DDE3F7)
$colorBlue = 16247774
$colorGray = 15921906
$colorLavender = 16245725

# --- 1. Remove the row describing the now-dropped "sends an update" event ---
# (was row 3: "Escalation Management Module sends an update of escalated ticket
# details to the escalation level employee ..."); everything below shifts up.
$ws.Rows("3:3").Delete()

# --- 2. Update wording on the remaining rows ---
$ws.Range("B2").Value = "Transforms service ticket ino an escalation ticket"
$ws.Range("D2").Value = "Receive Escalated Ticket"
$ws.Range("F2").Value = "Escalated Tickets Data Store "
$ws.Range("A4").Value = "Supervisor/Department Manager/Resident Manager/General Manager enter reason for escalated ticket "
$ws.Range("D4").Value = "Enter escalated ticket's reason"

# --- 3. Row heights to match the revised layout ---
$ws.Rows("5:5").RowHeight = 15
$ws.Rows("6:6").RowHeight = 89.25

# --- 4. Shading: header row ---
$ws.Range("A1:F1").Interior.Color = $colorBlue

# --- 5. Shading: banded data rows (2, 4, 6, 8) ---
$ws.Range("A2:F2").Interior.Color = $colorGray
$ws.Range("A4:F4").Interior.Color = $colorGray
$ws.Range("A6:F6").Interior.Color = $colorGray
$ws.Range("A8:F8").Interior.Color = $colorGray

# rows 3 and 7 keep the default (no fill)
$ws.Range("A3:F3").Interior.ColorIndex = -4142
$ws.Range("A7:F7").Interior.ColorIndex = -4142

# --- 6. Shading + formatting: the "MODULE" row (row 5) ---
$moduleRow = $ws.Range("A5:F5")
$moduleRow.Interior.Color = $colorLavender
$moduleRow.HorizontalAlignment = $xlCenter
$moduleRow.VerticalAlignment = $xlCenter
$moduleRow.WrapText = $true

# --- 7. Small alignment tweaks ---
$b7 = $ws.Range("B7")
$b7.HorizontalAlignment = $xlCenter
$b7.VerticalAlignment = $xlCenter
$b7.WrapText = $false

$b8 = $ws.Range("B8")
$b8.HorizontalAlignment = $xlCenter
$b8.VerticalAlignment = $xlCenter
$b8.WrapText = $false

# --- 8. Update selection to match the saved view ---
$ws.Range("A5").Select()

Write-Output "Event table revised."
